# Applies the "Actualización automática 2025-08-01 08:30:08" edit:
#  - Sheet "VENTAS POR GRUPO": zero out a handful of cells in rows 14/18/19
#    and update the "N de 27" progress labels in row 29.
#  - Sheet "VENTA MENSUAL": shift the reported month columns one position
#    to the right (abril->mayo->junio->julio->agosto) together with the
#    matching sales figures, and tweak two column widths.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Range("C14").Value = 0
$ws1.Range("L14").Value = 0
$ws1.Range("M14").Value = 0
$ws1.Range("N14").Value = 0
$ws1.Range("Q14").Value = 0

$ws1.Range("D18").Value = 0
$ws1.Range("M18").Value = 0

$ws1.Range("L19").Value = 0
$ws1.Range("M19").Value = 0

$ws1.Range("C29").Value = "0 de 27"
$ws1.Range("D29").Value = "0 de 27"
$ws1.Range("L29").Value = "0 de 27"
$ws1.Range("M29").Value = "0 de 27"

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Column widths (C and F). Excel's ColumnWidth property measures widths in
# characters and is offset from the raw OOXML "width" attribute by 5/6 of a
# character (default-font padding), so subtract that to land on the target
# stored widths of 13 and 12.
$ws2.Columns.Item(3).ColumnWidth = 13 - (5/6)
$ws2.Columns.Item(6).ColumnWidth = 12 - (5/6)

# Header row: month names shift one column to the right
$ws2.Range("C1").Value = "mayo"
$ws2.Range("D1").Value = "junio"
$ws2.Range("E1").Value = "julio"
$ws2.Range("F1").Value = "agosto"

# Row 4
$ws2.Range("C4").Value = 419.13
$ws2.Range("D4").Value = 0

# Row 13
$ws2.Range("D13").Value = 738.66
$ws2.Range("E13").Value = 0

# Row 14
$ws2.Range("C14").Value = 3122.02
$ws2.Range("D14").Value = 1473.73
$ws2.Range("E14").Value = 3990.41
$ws2.Range("F14").Value = 0

# Row 16
$ws2.Range("C16").Value = 0

# Row 17
$ws2.Range("D17").Value = 9556.26
$ws2.Range("E17").Value = 0

# Row 18
$ws2.Range("C18").Value = 6725.74
$ws2.Range("D18").Value = 64.81999999999999
$ws2.Range("E18").Value = 8691.84
$ws2.Range("F18").Value = 0

# Row 19
$ws2.Range("C19").Value = 0
$ws2.Range("D19").Value = 411.7
$ws2.Range("E19").Value = 2045.31
$ws2.Range("F19").Value = 0

# Row 21
$ws2.Range("C21").Value = 1994.73
$ws2.Range("D21").Value = 3225.33
$ws2.Range("E21").Value = 0

# Row 29 (totals)
$ws2.Range("C29").Value = 12261.62
$ws2.Range("D29").Value = 15470.5
$ws2.Range("E29").Value = 14727.56
$ws2.Range("F29").Value = 0
